$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 264.875
$ws.Range("I5").Value = 131.28572
$ws.Range("K5").Value = 131.28572
$ws.Range("M5").Value = -16.28572
$ws.Range("H18").Value = 743.7143
$ws.Range("I18").Value = 743.7143
$ws.Range("K18").Value = 743.7143
$ws.Range("M18").Value = -459.7143
$ws.Range("H92").Value = 50273.2
$ws.Range("I92").Value = 52813.895
$ws.Range("K92").Value = 52813.895
$ws.Range("M92").Value = -51565.895
$ws.Range("H118").Value = 975.13043
$ws.Range("I118").Value = 562.875
$ws.Range("K118").Value = 1688.625
$ws.Range("M118").Value = -31.625
$ws.Range("H138").Value = 3145.0625
$ws.Range("I138").Value = 1768.0714
$ws.Range("J138").Value = 3712.0588
$ws.Range("K138").Value = 5304.2142
$ws.Range("L138").Value = 11136.1764
$ws.Range("M138").Value = -164.2142000000003
$ws.Range("N138").Value = -21416.1764
$ws.Range("H141").Value = 10492.083
$ws.Range("I141").Value = 11163.637
$ws.Range("J141").Value = 3105
$ws.Range("K141").Value = 33490.911
$ws.Range("L141").Value = 9315
$ws.Range("M141").Value = -28310.911
$ws.Range("N141").Value = -19675

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 9566.091
$ws.Range("I32").Value = 5987.32
$ws.Range("K32").Value = 5987.32
$ws.Range("M32").Value = -5700.32
$ws.Range("H45").Value = 3784.9
$ws.Range("I45").Value = 1210.25
$ws.Range("K45").Value = 1210.25
$ws.Range("M45").Value = -833.25
$ws.Range("H61").Value = 3709.3635
$ws.Range("I61").Value = 2750.375
$ws.Range("K61").Value = 2750.375
$ws.Range("M61").Value = -2538.375
$ws.Range("H97").Value = 9650.5
$ws.Range("I97").Value = 15018.143
$ws.Range("K97").Value = 15018.143
$ws.Range("M97").Value = -14522.143
$ws.Range("H132").Value = 21254.02
$ws.Range("I132").Value = 22616.646
$ws.Range("K132").Value = 67849.93799999999
$ws.Range("M132").Value = -65319.93799999999
$ws.Range("H136").Value = 3709.3635
$ws.Range("I136").Value = 2750.375
$ws.Range("K136").Value = 8251.125
$ws.Range("M136").Value = -5701.125

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H75").Value = 9853.25
$ws.Range("I75").Value = 9853.25
$ws.Range("K75").Value = 9853.25
$ws.Range("M75").Value = -8917.25
$ws.Range("H78").Value = 9853.25
$ws.Range("I78").Value = 9853.25
$ws.Range("K78").Value = 29559.75
$ws.Range("M78").Value = -24879.75
$ws.Range("H82").Value = 27373.875
$ws.Range("I82").Value = 14666
$ws.Range("J82").Value = 65497.5
$ws.Range("K82").Value = 14666
$ws.Range("L82").Value = 65497.5
$ws.Range("M82").Value = -14283
$ws.Range("N82").Value = -66263.5
$ws.Range("H85").Value = 27373.875
$ws.Range("I85").Value = 14666
$ws.Range("J85").Value = 65497.5
$ws.Range("K85").Value = 14666
$ws.Range("L85").Value = 65497.5
$ws.Range("M85").Value = -13340
$ws.Range("N85").Value = -68149.5
$ws.Range("H86").Value = 2554.0908
$ws.Range("I86").Value = 2282.6667
$ws.Range("J86").Value = 2879.8
$ws.Range("K86").Value = 2282.6667
$ws.Range("L86").Value = 2879.8
$ws.Range("M86").Value = -1159.6667
$ws.Range("N86").Value = -5125.8
$ws.Range("H89").Value = 2554.0908
$ws.Range("I89").Value = 2282.6667
$ws.Range("J89").Value = 2879.8
$ws.Range("K89").Value = 11413.3335
$ws.Range("L89").Value = 14399
$ws.Range("M89").Value = -5797.333500000001
$ws.Range("N89").Value = -25631
$ws.Range("H116").Value = 74000
$ws.Range("J116").Value = 74000
$ws.Range("L116").Value = 74000
$ws.Range("N116").Value = -83178
$ws.Range("H134").Value = 1857.55
$ws.Range("I134").Value = 1640.6604
$ws.Range("J134").Value = 3499.7144
$ws.Range("K134").Value = 4921.9812
$ws.Range("L134").Value = 10499.1432
$ws.Range("M134").Value = -2386.9812
$ws.Range("N134").Value = -15569.1432

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 895.25
$ws.Range("I16").Value = 895.25
$ws.Range("K16").Value = 895.25
$ws.Range("M16").Value = -608.25
$ws.Range("H31").Value = 2448.0312
$ws.Range("I31").Value = 2052.9333
$ws.Range("J31").Value = 8374.5
$ws.Range("K31").Value = 2052.9333
$ws.Range("L31").Value = 8374.5
$ws.Range("M31").Value = -1757.9333
$ws.Range("N31").Value = -8964.5
$ws.Range("H34").Value = 2448.0312
$ws.Range("I34").Value = 2052.9333
$ws.Range("J34").Value = 8374.5
$ws.Range("K34").Value = 2052.9333
$ws.Range("L34").Value = 8374.5
$ws.Range("M34").Value = -1850.9333
$ws.Range("N34").Value = -8778.5
$ws.Range("H58").Value = 42228.56
$ws.Range("I58").Value = 45669.566
$ws.Range("K58").Value = 45669.566
$ws.Range("M58").Value = -45466.566
$ws.Range("H113").Value = 895.25
$ws.Range("I113").Value = 895.25
$ws.Range("K113").Value = 895.25
$ws.Range("M113").Value = 1274.75
$ws.Range("H132").Value = 2698.5898
$ws.Range("I132").Value = 2647.9707
$ws.Range("J132").Value = 3042.8
$ws.Range("K132").Value = 7943.9121
$ws.Range("L132").Value = 9128.400000000001
$ws.Range("M132").Value = -5413.9121
$ws.Range("N132").Value = -14188.4
$ws.Range("H136").Value = 42228.56
$ws.Range("I136").Value = 45669.566
$ws.Range("K136").Value = 137008.698
$ws.Range("M136").Value = -134458.698

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H37").Value = 68013.836
$ws.Range("J37").Value = 68013.836
$ws.Range("L37").Value = 204041.508
$ws.Range("N37").Value = -204265.508
$ws.Range("H113").Value = 1662.4667
$ws.Range("I113").Value = 1169.7142
$ws.Range("J113").Value = 1812.4348
$ws.Range("K113").Value = 3509.1426
$ws.Range("L113").Value = 5437.3044
$ws.Range("M113").Value = -1339.1426
$ws.Range("N113").Value = -9777.304400000001
$ws.Range("H121").Value = 741
$ws.Range("I121").Value = 340
$ws.Range("J121").Value = 923.2727
$ws.Range("K121").Value = 1020
$ws.Range("L121").Value = 2769.8181
$ws.Range("M121").Value = 290
$ws.Range("N121").Value = -5389.8181
$ws.Range("H122").Value = 491.13333
$ws.Range("I122").Value = 379.81818
$ws.Range("J122").Value = 797.25
$ws.Range("K122").Value = 3418.36362
$ws.Range("L122").Value = 7175.25
$ws.Range("M122").Value = -968.3636200000001
$ws.Range("N122").Value = -12075.25
$ws.Range("H132").Value = 1098.8
$ws.Range("I132").Value = 873.5
$ws.Range("K132").Value = 7861.5
$ws.Range("M132").Value = -5331.5
$ws.Range("H133").Value = 7374.8335
$ws.Range("J133").Value = 10062.5
$ws.Range("L133").Value = 30187.5
$ws.Range("N133").Value = -40307.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3647.2856
$ws.Range("I80").Value = 3200.3635
$ws.Range("K80").Value = 3200.3635
$ws.Range("M80").Value = -2202.3635
$ws.Range("H83").Value = 3647.2856
$ws.Range("I83").Value = 3200.3635
$ws.Range("K83").Value = 16001.8175
$ws.Range("M83").Value = -11009.8175
$ws.Range("H122").Value = 3799.5
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 3799.5
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 11398.5
$ws.Range("M122").ClearContents()
$ws.Range("N122").Value = -16298.5
$ws.Range("H132").Value = 29385.514
$ws.Range("J132").Value = 7429
$ws.Range("L132").Value = 22287
$ws.Range("N132").Value = -27347

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 67244.88
$ws.Range("I22").Value = 139420.88
$ws.Range("K22").Value = 139420.88
$ws.Range("M22").Value = -139125.88
$ws.Range("H27").Value = 67244.88
$ws.Range("I27").Value = 139420.88
$ws.Range("K27").Value = 139420.88
$ws.Range("M27").Value = -139313.88
$ws.Range("H46").Value = 15655.708
$ws.Range("J46").Value = 5341.1113
$ws.Range("L46").Value = 5341.1113
$ws.Range("N46").Value = -5717.1113
$ws.Range("H61").Value = 2741.4614
$ws.Range("I61").Value = 2251.12
$ws.Range("K61").Value = 2251.12
$ws.Range("M61").Value = -2049.12
$ws.Range("H113").Value = 2741.4614
$ws.Range("I113").Value = 2251.12
$ws.Range("K113").Value = 2251.12
$ws.Range("M113").Value = -81.11999999999989
$ws.Range("H136").Value = 2974.9048
$ws.Range("I136").Value = 2656.4736
$ws.Range("K136").Value = 7969.4208
$ws.Range("M136").Value = -5419.4208

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 9536
$ws.Range("I81").Value = 1881.5
$ws.Range("J81").Value = 32499.5
$ws.Range("K81").Value = 3763
$ws.Range("L81").Value = 64999
$ws.Range("M81").Value = -2702
$ws.Range("N81").Value = -67121
$ws.Range("H84").Value = 9536
$ws.Range("I84").Value = 1881.5
$ws.Range("J84").Value = 32499.5
$ws.Range("K84").Value = 18815
$ws.Range("L84").Value = 324995
$ws.Range("M84").Value = -13511
$ws.Range("N84").Value = -335603
$ws.Range("H107").Value = 1233.5555
$ws.Range("I107").Value = 1157.8572
$ws.Range("J107").Value = 1498.5
$ws.Range("K107").Value = 3473.5716
$ws.Range("L107").Value = 4495.5
$ws.Range("M107").Value = -1553.5716
$ws.Range("N107").Value = -8335.5
$ws.Range("H113").Value = 946.0357
$ws.Range("I113").Value = 417.2353
$ws.Range("J113").Value = 1763.2727
$ws.Range("K113").Value = 1251.7059
$ws.Range("L113").Value = 5289.8181
$ws.Range("M113").Value = 918.2941000000001
$ws.Range("N113").Value = -9629.8181
$ws.Range("H132").Value = 34437.355
$ws.Range("I132").Value = 36798.55
$ws.Range("K132").Value = 110395.65
$ws.Range("M132").Value = -107865.65
$ws.Range("H136").Value = 2411.3713
$ws.Range("I136").Value = 2031.7241
$ws.Range("K136").Value = 6095.1723
$ws.Range("M136").Value = -3545.1723
